# Apply updated crypto price/volume data (GitHub Actions refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").Value = "68.505.38"
$ws.Range("E2").Value = "  +1.63%  "

$ws.Range("D3").Value = "3.779.49"
$ws.Range("E3").Value = "  +0.44%  "

$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("E5").Value = "  +0.12%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "168.62"
$ws.Range("E6").Value = "  +0.12%  "

$ws.Range("D7").Value = "3.776.54"
$ws.Range("E7").Value = "  +0.40%  "

$ws.Range("E8").Value = "  -0.04%  "

$ws.Range("E9").Value = "  -1.27%  "

$ws.Range("E10").Value = "  -1.71%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.50"
$ws.Range("E11").Value = "  +0.25%  "

$ws.Range("E12").Value = "  -1.51%  "

$ws.Range("E13").Value = "  -3.19%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.53"
$ws.Range("E14").Value = "  -0.63%  "

$ws.Range("D15").Value = "4.413.71"
$ws.Range("E15").Value = "  +0.54%  "

$ws.Range("D16").Value = "3.776.98"
$ws.Range("E16").Value = "  +0.53%  "

$ws.Range("D17").Value = "68.473.07"
$ws.Range("E17").Value = "  +1.45%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "18.28"
$ws.Range("E18").Value = "  -3.75%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.06"
$ws.Range("E19").Value = "  -2.64%  "

$ws.Range("E20").Value = "  -0.30%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.00"
$ws.Range("E21").Value = "  +4.50%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "469.43"
$ws.Range("E22").Value = "  +0.43%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.703"

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "84.91"
$ws.Range("E24").Value = "  +1.31%  "

$ws.Range("E25").Value = "  -4.45%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.24"
$ws.Range("E26").Value = "  -0.04%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.19"
$ws.Range("E27").Value = "  +0.41%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.17"
$ws.Range("E28").Value = "  -1.11%  "

$ws.Range("E29").Value = "  +0.15%  "

$ws.Range("D30").Value = "3.926.23"
$ws.Range("E30").Value = "  +0.16%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.79"
$ws.Range("E31").Value = "  -3.69%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.41"
$ws.Range("E32").Value = "  -2.70%  "

$ws.Range("E33").Value = "  -1.12%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "30.11"
$ws.Range("E34").Value = "  -0.87%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.28"
$ws.Range("E35").Value = "  +1.49%  "

$ws.Range("D37").Value = "3.732.71"
$ws.Range("E37").Value = "  +0.24%  "

$ws.Range("E38").Value = "  -3.25%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.46"
$ws.Range("E39").Value = "  -9.72%  "

$ws.Range("E40").Value = "  +1.08%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.00"
$ws.Range("E41").Value = "  +0.28%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.82"
$ws.Range("E42").Value = "  -1.17%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.00"
$ws.Range("E43").Value = "  +0.02%  "

$ws.Range("E44").Value = "  +0.01%  "

$ws.Range("E45").Value = "  -2.32%  "

$ws.Range("E46").Value = "  +0.56%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "43.81"
$ws.Range("E47").Value = "  +12.66%  "

$ws.Range("E48").Value = "  -1.44%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "407.82"
$ws.Range("E49").Value = "  +1.56%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "45.69"
$ws.Range("E50").Value = "  -1.17%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "145.77"
$ws.Range("E51").Value = "  +2.67%  "
